$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column AT (46): "affix_type"
$ws.Cells.Item(1, 46).Value = "affix_type"

# Fill every data row (2-58) with the new affix_type value of 2
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 46).Value = 2
}

# Set the width of the new column to match the target (~12.65 chars)
$ws.Columns.Item(46).ColumnWidth = 11.8

# Update the active selection to the new column's data range
$null = $ws.Range("AT2:AT58").Select()
